# Auto-save via app Streamlit
# A new reservation (Gregory Blanvillain) was inserted as row 32 of the
# "reservations" sheet, pushing the previously-existing rows 32-47 down
# to rows 33-48 (the TOTAL row moves from row 47 to row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 32 ("Fermeture"), shifting
# every row below it (including the TOTAL row) down by one.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new reservation data.
$ws.Range("A32").Value = "Gregory Blanvillain"
$ws.Range("B32").Value = "Booking"

# Phone numbers are stored as text (leading "+"), so force text formatting
# before assigning the value, otherwise Excel would coerce it to a number.
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "+33687762155"

$ws.Range("D32").Value = 45877
$ws.Range("E32").Value = 45880
$ws.Range("F32").Value = 3
$ws.Range("G32").Value = 167.06
$ws.Range("H32").Value = 133.1
$ws.Range("I32").Value = 33.96
$ws.Range("J32").Value = 20.33
$ws.Range("K32").Value = 2025
$ws.Range("L32").Value = 8

# M32/N32/O32 stay blank. The row-insert operation copies the date
# formatting (style) from row 31 into N32/O32 - clear that back to the
# default (unstyled/blank) look used by the rest of the "Fermeture"-type
# rows.
$ws.Range("M32:O32").Style = "Normal"

# The "Fermeture" row that used to be row 32 is now row 33; its phone
# column held the placeholder text "nan" which should end up blank.
$ws.Range("C33").Value = ""
